$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PhanCongNganHan")

# Hoang added the "Quan ly thong tin Dia Diem" (3.6) results, plus filled
# in the two other results ("Quan ly thong tin Mat Hang" / row 10, and
# "Quan ly phan cong Xe - Tai Xe" / row 13) that reached 100% in this pass.
$ws.Range("E10").Value = "100% (11/06/2010)"
$ws.Range("E8").Value  = "100% (17/06/2010)"
$ws.Range("E13").Value = "100% (16/06/2010)"

# Leave the view focused where the author was last working.
$ws.Activate()
$ws.Range("E13").Select()
